# Pomc-Mc1r.xlsx - "update scripts wuth new tpm"
#
# NATMI LR-pair output recomputed with updated TPM values:
#   - existing "ECs -> *" rows (2-6) get refreshed edge/ligand/receptor stats
#   - a new sending cluster "MuSCs" is introduced, contributing 5 more rows
#     (7-11) against the same ordered set of target clusters
#     (ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-T per row: SendingCluster, Ligand, Receptor, TargetCluster,
# then the 16 numeric NATMI statistics columns (E..T).
$rows = @(
   ,@(2,  "ECs",   "Pomc", "Mc1r", "ECs",               3, 1,                  0.5131323333333334, 1.539397, 0.9119558630037493,  0.9119558630037494,  3, 1, 0.6201396666666666, 1.860419, 0.1891028895630684, 0.1891028895630684, 0.3182137141492223,  2.863923427343,  0.1724534888479908,  0.1724534888479908)
   ,@(3,  "ECs",   "Pomc", "Mc1r", "FAPs",              3, 1,                  0.5131323333333334, 1.539397, 0.9119558630037493,  0.9119558630037494,  3, 1, 1.018725333333333,  3.056176, 0.3106459956672665, 0.3106459956672664, 0.5227409073191112,  4.704668165872,  0.283295437067401,   0.283295437067401)
   ,@(4,  "ECs",   "Pomc", "Mc1r", "Inflammatory-Mac",  3, 1,                  0.5131323333333334, 1.539397, 0.9119558630037493,  0.9119558630037494,  3, 1, 0.5912873333333333, 1.773862, 0.1803047753684109, 0.1803047753684109, 0.3034086490237778,  2.730677841214,  0.1644299970247963,  0.1644299970247964)
   ,@(5,  "ECs",   "Pomc", "Mc1r", "MuSCs",             3, 1,                  0.5131323333333334, 1.539397, 0.9119558630037493,  0.9119558630037494,  3, 1, 0.4518883333333334, 1.355665, 0.1377970063622857, 0.1377970063622857, 0.2318785148894445,  2.086906634005,  0.1256647878564514,  0.1256647878564514)
   ,@(6,  "ECs",   "Pomc", "Mc1r", "Resolving-Mac",     3, 1,                  0.5131323333333334, 1.539397, 0.9119558630037493,  0.9119558630037494,  3, 1, 0.5973363333333334, 1.792009, 0.1821493330389685, 0.1821493330389685, 0.3065125865081112,  2.758613278573,  0.1661121522071099,  0.1661121522071098)
   ,@(7,  "MuSCs", "Pomc", "Mc1r", "ECs",               1, 0.3333333333333333, 0.04954,             0.14862,  0.08804413699625062,  0.08804413699625063, 3, 1, 0.6201396666666666, 1.860419, 0.1891028895630684, 0.1891028895630684, 0.03072171908666666, 0.27649547178,   0.01664940071507765, 0.01664940071507765)
   ,@(8,  "MuSCs", "Pomc", "Mc1r", "FAPs",              1, 0.3333333333333333, 0.04954,             0.14862,  0.08804413699625062,  0.08804413699625063, 3, 1, 1.018725333333333,  3.056176, 0.3106459956672665, 0.3106459956672664, 0.05046765301333333, 0.45420887712,   0.02735055859986548, 0.02735055859986548)
   ,@(9,  "MuSCs", "Pomc", "Mc1r", "Inflammatory-Mac",  1, 0.3333333333333333, 0.04954,             0.14862,  0.08804413699625062,  0.08804413699625063, 3, 1, 0.5912873333333333, 1.773862, 0.1803047753684109, 0.1803047753684109, 0.02929237449333333, 0.26363137044,   0.01587477834361456, 0.01587477834361457)
   ,@(10, "MuSCs", "Pomc", "Mc1r", "MuSCs",             1, 0.3333333333333333, 0.04954,             0.14862,  0.08804413699625062,  0.08804413699625063, 3, 1, 0.4518883333333334, 1.355665, 0.1377970063622857, 0.1377970063622857, 0.02238654803333334, 0.2014789323,    0.0121322185058343,  0.0121322185058343)
   ,@(11, "MuSCs", "Pomc", "Mc1r", "Resolving-Mac",     1, 0.3333333333333333, 0.04954,             0.14862,  0.08804413699625062,  0.08804413699625063, 3, 1, 0.5973363333333334, 1.792009, 0.1821493330389685, 0.1821493330389685, 0.02959204195333334, 0.26632837758,   0.01603718083185862, 0.01603718083185862)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $r[$col]
    }
}
